# Generate Report for Handback
# Refresh the "latest generated/handoff/handback" timestamps written by the
# handback status report generator.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-14 03:36:02"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-14 03:35:53"
$wsZhCn.Range("K2").Value = "2016-08-14 03:36:17"

# de-de sheet: Correspond Handoff Datetime
$wsDeDe.Range("H2").Value = "2016-08-14 03:36:28"
